$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Objetivos:" / "Objectives:" gains real body text -----------
# (previously B10/C10 incorrectly held the professor's name)
$ws.Cells.Item(10, 2).Value = "Apresentar conceitos, ferramentas e métodos para o auxílio à tomada de decisão."
$ws.Cells.Item(10, 3).Value = "Apresentar conceitos, ferramentas e métodos para o auxílio à tomada de decisão."

# --- Insert a new row at 13 to hold the professor's name on its own row --
$ws.Rows.Item(13).Insert()

# The inserted row drags column A's formatting down into the new row 13;
# the target layout has no A13 cell at all, so clear it...
$ws.Cells.Item(13, 1).Clear()

# ...then borrow B/C number formats from row 14 (same columns, still
# holding the old "Programa resumido" body at this point) so the new
# B13/C13 cells keep the normal body style (s="2"/s="3") instead of
# falling back to column A's style.
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B13:C13").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(13, 2).Value = "5840917 - Fabricio Maciel Gomes"
$ws.Cells.Item(13, 3).Value = "5840917 - Fabricio Maciel Gomes"

# --- Row 14 (old row 13): "Programa resumido:" body text ------------------
$ws.Cells.Item(14, 2).Value = "Teoria da Decisão; Estruturação, Decisão sem Risco e sem Incerteza; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão e Sistemas Especialistas."
$ws.Cells.Item(14, 3).Value = "Teoria da Decisão; Estruturação, Decisão sem Risco e sem Incerteza; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão e Sistemas Especialistas."

# --- Row 16 (old row 15): "Programa:" body text ----------------------------
$ws.Cells.Item(16, 2).Value = "1.Teoria da Decisão`n2.Estruturação, Decisão sem Risco e sem Incerteza`n3.Decisão com Múltiplos Cenários ou Múltiplos Critérios`n4.Decisão com Incerteza`n5.Sistemas de Auxílio à Decisão e Sistemas Especialistas."
$ws.Cells.Item(16, 3).Value = "1.Teoria da Decisão`n2.Estruturação, Decisão sem Risco e sem Incerteza`n3.Decisão com Múltiplos Cenários ou Múltiplos Critérios`n4.Decisão com Incerteza`n5.Sistemas de Auxílio à Decisão e Sistemas Especialistas."

# --- Row 19 (old row 18): "Método:" body text ------------------------------
$ws.Cells.Item(19, 2).Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Cells.Item(19, 3).Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."

# --- Row 20 (old row 19): "Critério:" body text ----------------------------
$ws.Cells.Item(20, 2).Value = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"
$ws.Cells.Item(20, 3).Value = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"

# --- Row 21 (old row 20): "Norma de recuperação:" body text ---------------
$ws.Cells.Item(21, 2).Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Cells.Item(21, 3).Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."

# --- Row 22 (old row 21): "Bibliografia:" body text ------------------------
$ws.Cells.Item(22, 2).Value = "1.ENSSLIN, L.; MONTIBELLER NETO, G.; NORONHA, S. M. (2001), Apoio à Decisão: metodologias para estruturação de problemas e avaliação multicritério de alternativas, Florianópolis: Insular`n2.GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T. (2002), Tomada de Decisão Gerencial: enfoque multicritério, São Paulo: Atlas`n3.LAWRENCE, J. A. JR.; PASTERNACK, B. A. (2002), Applied Management Science: modeling, spreadsheet analysis and communication for decision making, 2nd edition, New York (USA): Wiley`n4.SHIMIZU, T. (2001), Decisão nas Organizações: introdução aos problemas de decisão encontrados nas organizações e nos sistemas de apoio à decisão, São Paulo: Atlas"
$ws.Cells.Item(22, 3).Value = "1.ENSSLIN, L.; MONTIBELLER NETO, G.; NORONHA, S. M. (2001), Apoio à Decisão: metodologias para estruturação de problemas e avaliação multicritério de alternativas, Florianópolis: Insular`n2.GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T. (2002), Tomada de Decisão Gerencial: enfoque multicritério, São Paulo: Atlas`n3.LAWRENCE, J. A. JR.; PASTERNACK, B. A. (2002), Applied Management Science: modeling, spreadsheet analysis and communication for decision making, 2nd edition, New York (USA): Wiley`n4.SHIMIZU, T. (2001), Decisão nas Organizações: introdução aos problemas de decisão encontrados nas organizações e nos sistemas de apoio à decisão, São Paulo: Atlas"
